$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: update title (D23) and link (E23)
$ws.Range("D23").Value = "Google의 머신러닝 엔지니어링 실무 지침서입니다.`n머신러닝 프로젝트 구조화에 대해 공부하다가 찾게 된 문서인데, 내용이 너무 좋아 공유 드"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2691"

# Row 28: update title (D28) and link (E28)
$ws.Range("D28").Value = "DRQN 구현"
$ws.Range("E28").Value = "https://ropiens.tistory.com/80"

# Row 46: update title (D46) and link (E46)
$ws.Range("D46").Value = "심전도의 3가지 규칙"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/372"
